# "update review with PGx team"
#
# The PGx team re-reviewed sample 20220112-24001 and corrected a couple of
# bin boundaries on the peak_table sheet, which flips which CYP2D6_011 /
# CYP2D6_012 allele calls are considered "detected" on the allele_table
# sheet, which in turn flips the genotype calls on the marker_table sheet
# and fills in the final diplotype on the genotype_result sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# peak_table: corrected bin boundaries (highlighted in red so the PGx
# team can see exactly what changed during their review)
# ---------------------------------------------------------------------
$peak = $wb.Worksheets.Item("peak_table")

$peak.Range("G12").Value = 33
$peak.Range("G12").Font.Color = 255

$peak.Range("H13").Value = 33
$peak.Range("H13").Font.Color = 255

$peak.Activate()
$peak.Range("H13").Select()

# ---------------------------------------------------------------------
# allele_table: with the corrected bins, CYP2D6_011 (row 22) now detects
# and CYP2D6_012 (row 25) no longer does - swap the detection results.
# ---------------------------------------------------------------------
$allele = $wb.Worksheets.Item("allele_table")

$allele.Range("J22").Value = 33
$allele.Range("M22").Value = $true
$allele.Range("N22").Value = 58
$allele.Range("O22").Value = 32.2
$allele.Range("P22").Value = 2131
$allele.Range("Q22").Value = "ok"
$allele.Range("R22").Value = ""

$allele.Range("I25").Value = 33
$allele.Range("M25").Value = $false
$allele.Range("N25").Value = ""
$allele.Range("O25").Value = ""
$allele.Range("P25").Value = ""
$allele.Range("Q25").Value = ""
$allele.Range("R25").Value = "Peak(s) could not be detected. Please check peak ranges if required!"

# ---------------------------------------------------------------------
# marker_table: genotype calls flip accordingly
# ---------------------------------------------------------------------
$marker = $wb.Worksheets.Item("marker_table")

$marker.Range("G12").Value = "GA"
$marker.Range("H12").Value = "heterozygous"

$marker.Range("G13").Value = "AA"
$marker.Range("H13").Value = "wildtype"

# ---------------------------------------------------------------------
# genotype_result: final reviewed diplotype call
# ---------------------------------------------------------------------
$result = $wb.Worksheets.Item("genotype_result")
$result.Range("B2").Value = "*2/*4"

# ---------------------------------------------------------------------
# Header styling refresh applied by the PGx team across the data sheets
# (bold header row with a thin box border, centered/top aligned).
# ---------------------------------------------------------------------
foreach ($sheetName in @("allele_table", "marker_table", "genotype_result")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $ws.UsedRange.Columns.Count))
    $headerRange.Font.Bold = $true
    $headerRange.Borders.LineStyle = 1
    $headerRange.Borders.Weight = 2
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
}

$wb.Save()
